# AP20_TestData_Manage Transactions - Standard Invoices_21C.xlsx
# "Add files via upload" / "Anu - AP Files Uploaded"
#
# The Input_Value sheet had its Y2:AA2 "helper" cells (login URL / username /
# password, duplicated from Sheet1) wiped out and the hyperlink on Y2 removed
# when the file was re-uploaded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Select the range that is about to be edited (matches the author having
# highlighted Y2:AA2 before deleting its contents).
$ws.Range("Y2:AA2").Select()

# Drop the hyperlink that lived on Y2 (https://edrx.fa.us2.oraclecloud.com/).
$ws.Range("Y2").Hyperlinks.Delete()

# Clear the values of Y2, Z2 and AA2 - the cell formatting/styles stay as-is.
$ws.Range("Y2:AA2").ClearContents()
